$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = "Computer Science"
$ws.Range("C17").Value = 1.051825159337061
$ws.Range("B18").Value = "AI"
$ws.Range("C18").Value = 1.309345107064612
$ws.Range("B22").Value = "Computer Science"
$ws.Range("C22").Value = 1.061984453095498
$ws.Range("B23").Value = "AI"
$ws.Range("C23").Value = 1.32078453497408
$ws.Range("B62").Value = "Computer Science"
$ws.Range("C62").Value = 1.149662857102538
$ws.Range("B63").Value = "AI"
$ws.Range("C63").Value = 0.8973516234105906
$ws.Range("B66").Value = "Computer Science"
$ws.Range("C66").Value = 1.058689616851366
$ws.Range("B67").Value = "AI"
$ws.Range("C67").Value = 0.8117268046484413
$ws.Range("B68").Value = "Levels"
$ws.Range("C68").Value = 0.8023308361177167
$ws.Range("B82").Value = "Computer Science"
$ws.Range("C82").Value = 1.052679628979815
$ws.Range("B83").Value = "AI"
$ws.Range("C83").Value = 1.310307363760873
$ws.Range("B91").Value = "Computer Science"
$ws.Range("C91").Value = 0.9354285654204095
$ws.Range("B92").Value = "Levels"
$ws.Range("C92").Value = 1.177096843518948
$ws.Range("B95").Value = "Finance"
$ws.Range("C95").Value = 1.004635220030704
$ws.Range("B96").Value = "Electricity"
$ws.Range("C96").Value = 0.2885249267089345
$ws.Range("B102").Value = "Computer Science"
$ws.Range("C102").Value = 1.046991585851753
$ws.Range("B103").Value = "AI"
$ws.Range("C103").Value = 1.084087129779844
$ws.Range("B107").Value = "Banking"
$ws.Range("C107").Value = 0.402040483774579
$ws.Range("B108").Value = "Computer Science"
$ws.Range("C108").Value = 1.164768315389567
$ws.Range("B109").Value = "Cloud Computing"
$ws.Range("C109").Value = 0.003015071025280447
$ws.Range("B124").Value = "Computer Science"
$ws.Range("C124").Value = 1.045324440417545
$ws.Range("B125").Value = "AI"
$ws.Range("C125").Value = 1.082389595924046
$ws.Range("B131").Value = "Finance"
$ws.Range("C131").Value = 0.9977525088168081
$ws.Range("B132").Value = "AI"
$ws.Range("C132").Value = 1.312131528785209
$ws.Range("B140").Value = "Computer Science"
$ws.Range("C140").Value = 1.046142100548664
$ws.Range("B141").Value = "AI"
$ws.Range("C141").Value = 1.083222161725817
$ws.Range("B146").Value = "Computer Science"
$ws.Range("C146").Value = 1.052032734757469
$ws.Range("B147").Value = "AI"
$ws.Range("C147").Value = 1.309578869092147
$ws.Range("B156").Value = "Computer Science"
$ws.Range("C156").Value = 1.045867975200999
$ws.Range("B157").Value = "AI"
$ws.Range("C157").Value = 1.082943039489161
$ws.Range("B159").Value = "Finance"
$ws.Range("C159").Value = 1.005197686377464
$ws.Range("B160").Value = "Electricity"
$ws.Range("C160").Value = 0.2888858810129574
$ws.Range("B164").Value = "Computer Science"
$ws.Range("C164").Value = 1.153008121267897
$ws.Range("B165").Value = "AI"
$ws.Range("C165").Value = 0.9003023474705429
$ws.Range("B166").Value = "Computer Science"
$ws.Range("C166").Value = 1.034777040740063
$ws.Range("B167").Value = "AI"
$ws.Range("C167").Value = 1.290141852419778
$ws.Range("B168").Value = "Levels"
$ws.Range("C168").Value = 1.08734780026408
$ws.Range("B170").Value = "Computer Science"
$ws.Range("C170").Value = 1.049620302549755
$ws.Range("B171").Value = "Levels"
$ws.Range("C171").Value = 1.102577451388186
$ws.Range("B181").Value = "Computer Science"
$ws.Range("C181").Value = 1.045402101449318
$ws.Range("B182").Value = "AI"
$ws.Range("C182").Value = 1.082468672815592
$ws.Range("B185").Value = "Computer Science"
$ws.Range("C185").Value = 1.054868618009125
$ws.Range("B186").Value = "AI"
$ws.Range("C186").Value = 0.8082442268551389
$ws.Range("B187").Value = "Levels"
$ws.Range("C187").Value = 0.7988463550005407
$ws.Range("B205").Value = "Computer Science"
$ws.Range("C205").Value = 1.039247667390604
$ws.Range("B206").Value = "DevOps"
$ws.Range("C206").Value = 1.080289987075832
$ws.Range("B214").Value = "Computer Science"
$ws.Range("C214").Value = 1.047030124807407
$ws.Range("B215").Value = "Levels"
$ws.Range("C215").Value = 1.294946667666711
$ws.Range("B279").Value = "Computer Science"
$ws.Range("C279").Value = 1.049237009245282
$ws.Range("B280").Value = "Hardware"
$ws.Range("C280").Value = 0.6905856838873534
$ws.Range("B331").Value = "Finance"
$ws.Range("C331").Value = 0.9965912285912136
$ws.Range("B332").Value = "AI"
$ws.Range("C332").Value = 1.086979311609661
$ws.Range("B377").Value = "Computer Science"
$ws.Range("C377").Value = 1.046277668348691
$ws.Range("B378").Value = "Levels"
$ws.Range("C378").Value = 1.099147923823703
$ws.Range("B380").Value = "Computer Science"
$ws.Range("C380").Value = 1.064043360835601
$ws.Range("B381").Value = "Finance"
$ws.Range("C381").Value = 1.010432249814712
$ws.Range("B382").Value = "Computer Science"
$ws.Range("C382").Value = 1.056615399094345
$ws.Range("B383").Value = "AI"
$ws.Range("C383").Value = 1.314739329745646
$ws.Range("B384").Value = "Levels"
$ws.Range("C384").Value = 1.109754165537558
$ws.Range("B401").Value = "Computer Science"
$ws.Range("C401").Value = 1.049407807327577
$ws.Range("B402").Value = "AI"
$ws.Range("C402").Value = 1.306622699698165
$ws.Range("B403").Value = "Healthcare"
$ws.Range("C403").Value = 1.085710095892315
$ws.Range("B409").Value = "Levels"
$ws.Range("C409").Value = 1.12265500806447
$ws.Range("B410").Value = "Computer Science"
$ws.Range("C410").Value = 1.069190500093204
$ws.Range("B411").Value = "AI"
$ws.Range("C411").Value = 1.328896734761429
$ws.Range("B412").Value = "Hardware"
$ws.Range("C412").Value = 1.383283339117241
$ws.Range("B420").Value = "Computer Science"
$ws.Range("C420").Value = 1.064232031979684
$ws.Range("B421").Value = "DevOps"
$ws.Range("C421").Value = 0.4093440614403802
$ws.Range("B428").Value = "Computer Science"
$ws.Range("C428").Value = 0.8604093657647406
$ws.Range("B429").Value = "Levels"
$ws.Range("C429").Value = 1.095362611778106
$ws.Range("B439").Value = "Computer Science"
$ws.Range("C439").Value = 1.065225681702834
$ws.Range("B440").Value = "AI"
$ws.Range("C440").Value = 1.102652809700569
$ws.Range("B441").Value = "Levels"
$ws.Range("C441").Value = 1.315353527948743
$ws.Range("B453").Value = "Computer Science"
$ws.Range("C453").Value = 1.166975145422432
$ws.Range("B454").Value = "Levels"
$ws.Range("C454").Value = 0.9030071338759323
$ws.Range("B456").Value = "Computer Science"
$ws.Range("C456").Value = 1.171380779405432
$ws.Range("B458").Value = "Cybersecurity"
$ws.Range("C458").Value = 0
$ws.Range("B474").Value = "Finance"
$ws.Range("C474").Value = 1.26867825078507
$ws.Range("B475").Value = "AI"
$ws.Range("C475").Value = 0.8811088670264807
